$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 172, shifting existing rows 172:280 down to 173:281
$ws.Rows.Item(172).Insert()

# Populate the newly inserted row 172 with the new weekly record
$ws.Cells.Item(172, 1).Value = 8
$ws.Cells.Item(172, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(172, 3).Value = "Coquimbo"
$ws.Cells.Item(172, 4).Value = 44596
$ws.Cells.Item(172, 5).Value = 4
$ws.Cells.Item(172, 6).Value = 100114013
$ws.Cells.Item(172, 7).Value = "Zanahoria"
$ws.Cells.Item(172, 8).Value = "Sin especificar"
$ws.Cells.Item(172, 9).Value = "Primera"
$ws.Cells.Item(172, 10).Value = 760
$ws.Cells.Item(172, 11).Value = 5500
$ws.Cells.Item(172, 12).Value = 6000
$ws.Cells.Item(172, 13).Value = 5750
$ws.Cells.Item(172, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(172, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(172, 16).Value = 288
$ws.Cells.Item(172, 17).Value = 20
$ws.Cells.Item(172, 18).Value = "Hortaliza"
